# Append four more match rows (6-9) to the "Isuru Udana " batting log.
# These are the same match rows already present (rows 2-5), re-listed
# further down the sheet - row 6 repeats row 3, row 7 repeats row 2,
# row 8 repeats row 5 and row 9 repeats row 4.
#
# Using Copy/PasteSpecial(xlPasteValues) instead of Range.Value so that
# numeric-looking text (e.g. "4", "200.00") is preserved as TEXT, exactly
# like the source cells, instead of being coerced into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$ws.Range("A3:K3").Copy()
$ws.Range("A6:K6").PasteSpecial($xlPasteValues)

$ws.Range("A2:K2").Copy()
$ws.Range("A7:K7").PasteSpecial($xlPasteValues)

$ws.Range("A5:K5").Copy()
$ws.Range("A8:K8").PasteSpecial($xlPasteValues)

$ws.Range("A4:K4").Copy()
$ws.Range("A9:K9").PasteSpecial($xlPasteValues)
